$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3510
$ws.Range("I64").Value = 3033.3333
$ws.Range("J64").Value = 4225
$ws.Range("K64").Value = 3033.3333
$ws.Range("L64").Value = 4225
$ws.Range("M64").Value = -2785.3333
$ws.Range("N64").Value = -4721
$ws.Range("H67").Value = 3510
$ws.Range("I67").Value = 3033.3333
$ws.Range("J67").Value = 4225
$ws.Range("K67").Value = 3033.3333
$ws.Range("L67").Value = 4225
$ws.Range("M67").Value = -2175.3333
$ws.Range("N67").Value = -5941
$ws.Range("H70").Value = 2892.5
$ws.Range("I70").Value = 1520
$ws.Range("J70").Value = 3872.8572
$ws.Range("K70").Value = 4560
$ws.Range("L70").Value = 11618.5716
$ws.Range("M70").Value = -4290
$ws.Range("N70").Value = -12158.5716
$ws.Range("H73").Value = 2892.5
$ws.Range("I73").Value = 1520
$ws.Range("J73").Value = 3872.8572
$ws.Range("K73").Value = 4560
$ws.Range("L73").Value = 11618.5716
$ws.Range("M73").Value = -3624
$ws.Range("N73").Value = -13490.5716
$ws.Range("H112").Value = 1285
$ws.Range("J112").Value = 1322.9546
$ws.Range("L112").Value = 3968.8638
$ws.Range("N112").Value = -6184.8638
$ws.Range("H121").Value = 502.6842
$ws.Range("J121").Value = 502.6842
$ws.Range("L121").Value = 1508.0526
$ws.Range("N121").Value = -5002.0526
$ws.Range("H130").Value = 41890
$ws.Range("J130").Value = 41890
$ws.Range("L130").Value = 41890
$ws.Range("N130").Value = -51930
$ws.Range("H132").Value = 29147796
$ws.Range("I132").Value = 32585660
$ws.Range("K132").Value = 97756980
$ws.Range("M132").Value = -97754450
$ws.Range("H137").Value = 605466.5
$ws.Range("I137").Value = 1445765.9
$ws.Range("J137").Value = 2643.0435
$ws.Range("K137").Value = 4337297.699999999
$ws.Range("L137").Value = 7929.130500000001
$ws.Range("M137").Value = -4334747.699999999
$ws.Range("N137").Value = -13029.1305
$ws.Range("H138").Value = 2529.9644
$ws.Range("I138").Value = 1643.4706
$ws.Range("J138").Value = 3900
$ws.Range("K138").Value = 4930.4118
$ws.Range("L138").Value = 11700
$ws.Range("M138").Value = 209.5882000000001
$ws.Range("N138").Value = -21980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5085.365
$ws.Range("I32").Value = 4913.8076
$ws.Range("J32").Value = 5896.364
$ws.Range("K32").Value = 4913.8076
$ws.Range("L32").Value = 5896.364
$ws.Range("M32").Value = -4626.8076
$ws.Range("N32").Value = -6470.364
$ws.Range("H41").Value = 13120.071
$ws.Range("I41").Value = 1342
$ws.Range("J41").Value = 28824.166
$ws.Range("K41").Value = 1342
$ws.Range("L41").Value = 28824.166
$ws.Range("M41").Value = -928
$ws.Range("N41").Value = -29652.166
$ws.Range("H76").Value = 38700
$ws.Range("J76").Value = 38700
$ws.Range("L76").Value = 38700
$ws.Range("N76").Value = -39376
$ws.Range("H79").Value = 38700
$ws.Range("J79").Value = 38700
$ws.Range("L79").Value = 38700
$ws.Range("N79").Value = -41040
$ws.Range("H102").Value = 2066
$ws.Range("I102").Value = 2066
$ws.Range("K102").Value = 2066
$ws.Range("M102").Value = -444
$ws.Range("H125").Value = 41706.363
$ws.Range("J125").Value = 41706.363
$ws.Range("L125").Value = 41706.363
$ws.Range("N125").Value = -51546.363
$ws.Range("H132").Value = 2537.608
$ws.Range("I132").Value = 2046.425
$ws.Range("J132").Value = 4323.727
$ws.Range("K132").Value = 6139.275
$ws.Range("L132").Value = 12971.181
$ws.Range("M132").Value = -3609.275
$ws.Range("N132").Value = -18031.181
$ws.Range("H137").Value = 38141.332
$ws.Range("J137").Value = 38141.332
$ws.Range("L137").Value = 38141.332
$ws.Range("N137").Value = -48341.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 500.8
$ws.Range("I94").Value = 444.22223
$ws.Range("J94").Value = 1010
$ws.Range("K94").Value = 444.22223
$ws.Range("L94").Value = 1010
$ws.Range("M94").Value = 6.777769999999975
$ws.Range("N94").Value = -1912
$ws.Range("H99").Value = 2626.3572
$ws.Range("I99").Value = 1332.1111
$ws.Range("J99").Value = 4956
$ws.Range("K99").Value = 1332.1111
$ws.Range("L99").Value = 4956
$ws.Range("M99").Value = 165.8888999999999
$ws.Range("N99").Value = -7952
$ws.Range("H104").Value = 49800
$ws.Range("J104").Value = 49800
$ws.Range("L104").Value = 49800
$ws.Range("N104").Value = -56788

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 112.29412
$ws.Range("I38").Value = 59.9
$ws.Range("J38").Value = 187.14285
$ws.Range("K38").Value = 179.7
$ws.Range("L38").Value = 561.4285500000001
$ws.Range("M38").Value = 167.3
$ws.Range("N38").Value = -1255.42855
$ws.Range("H57").Value = 3005
$ws.Range("I57").Value = 3005
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 9015
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -8456
$ws.Range("N57").ClearContents()
$ws.Range("H68").Value = 2971.85
$ws.Range("I68").Value = 1108.7826
$ws.Range("J68").Value = 4129.973
$ws.Range("K68").Value = 3326.3478
$ws.Range("L68").Value = 12389.919
$ws.Range("M68").Value = -2515.3478
$ws.Range("N68").Value = -14011.919
$ws.Range("H71").Value = 2971.85
$ws.Range("I71").Value = 1108.7826
$ws.Range("J71").Value = 4129.973
$ws.Range("K71").Value = 9979.0434
$ws.Range("L71").Value = 37169.757
$ws.Range("M71").Value = -5923.0434
$ws.Range("N71").Value = -45281.757
$ws.Range("H80").Value = 8633.941000000001
$ws.Range("I80").Value = 7666.6665
$ws.Range("J80").Value = 8841.214
$ws.Range("K80").Value = 22999.9995
$ws.Range("L80").Value = 26523.642
$ws.Range("M80").Value = -22063.9995
$ws.Range("N80").Value = -28395.642
$ws.Range("H83").Value = 8633.941000000001
$ws.Range("I83").Value = 7666.6665
$ws.Range("J83").Value = 8841.214
$ws.Range("K83").Value = 68999.9985
$ws.Range("L83").Value = 79570.92600000001
$ws.Range("M83").Value = -64319.9985
$ws.Range("N83").Value = -88930.92600000001
$ws.Range("H86").Value = 897.913
$ws.Range("I86").Value = 457.69232
$ws.Range("J86").Value = 1470.2
$ws.Range("K86").Value = 1373.07696
$ws.Range("L86").Value = 4410.6
$ws.Range("M86").Value = -187.0769599999999
$ws.Range("N86").Value = -6782.6
$ws.Range("H89").Value = 897.913
$ws.Range("I89").Value = 457.69232
$ws.Range("J89").Value = 1470.2
$ws.Range("K89").Value = 4119.23088
$ws.Range("L89").Value = 13231.8
$ws.Range("M89").Value = 1808.76912
$ws.Range("N89").Value = -25087.8
$ws.Range("H113").Value = 1866273.8
$ws.Range("I113").Value = 606.6111
$ws.Range("J113").Value = 9615968
$ws.Range("K113").Value = 1819.8333
$ws.Range("L113").Value = 28847904
$ws.Range("M113").Value = 350.1667000000002
$ws.Range("N113").Value = -28852244
$ws.Range("H131").Value = 791.0526
$ws.Range("I131").Value = 440.36365
$ws.Range("J131").Value = 836.9761999999999
$ws.Range("K131").Value = 1321.09095
$ws.Range("L131").Value = 2510.9286
$ws.Range("M131").Value = 3718.90905
$ws.Range("N131").Value = -12590.9286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 35000
$ws.Range("J27").Value = 35000
$ws.Range("L27").Value = 35000
$ws.Range("N27").Value = -35332
$ws.Range("H70").Value = 6181.156
$ws.Range("I70").Value = 5514.676
$ws.Range("J70").Value = 9263.625
$ws.Range("K70").Value = 5514.676
$ws.Range("L70").Value = 9263.625
$ws.Range("M70").Value = -5244.676
$ws.Range("N70").Value = -9803.625
$ws.Range("H73").Value = 6181.156
$ws.Range("I73").Value = 5514.676
$ws.Range("J73").Value = 9263.625
$ws.Range("K73").Value = 5514.676
$ws.Range("L73").Value = 9263.625
$ws.Range("M73").Value = -4578.676
$ws.Range("N73").Value = -11135.625
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2658.1
$ws.Range("I132").Value = 2023.3096
$ws.Range("K132").Value = 6069.9288
$ws.Range("M132").Value = -3539.9288
$ws.Range("H133").Value = 63085
$ws.Range("J133").Value = 63085
$ws.Range("L133").Value = 63085
$ws.Range("N133").Value = -73205

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 102170
$ws.Range("I22").Value = 143814.28
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 143814.28
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -143519.28
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 102170
$ws.Range("I27").Value = 143814.28
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 143814.28
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -143707.28
$ws.Range("N27").Value = -5214
$ws.Range("H93").Value = 4833054
$ws.Range("I93").Value = 6537820
$ws.Range("J93").Value = 2884.6667
$ws.Range("K93").Value = 6537820
$ws.Range("L93").Value = 2884.6667
$ws.Range("M93").Value = -6536572
$ws.Range("N93").Value = -5380.6667
$ws.Range("H132").Value = 3833.3447
$ws.Range("I132").Value = 2757.5
$ws.Range("J132").Value = 7214.5713
$ws.Range("K132").Value = 8272.5
$ws.Range("L132").Value = 21643.7139
$ws.Range("M132").Value = -5742.5
$ws.Range("N132").Value = -26703.7139
$ws.Range("H136").Value = 4886.16
$ws.Range("I136").Value = 2877.182
$ws.Range("J136").Value = 6464.643
$ws.Range("K136").Value = 8631.545999999998
$ws.Range("L136").Value = 19393.929
$ws.Range("M136").Value = -6081.545999999998
$ws.Range("N136").Value = -24493.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 34970.332
$ws.Range("I29").Value = 20000
$ws.Range("J29").Value = 42455.5
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 42455.5
$ws.Range("M29").Value = -19710
$ws.Range("N29").Value = -43035.5
$ws.Range("H131").Value = 41640
$ws.Range("J131").Value = 41913.332
$ws.Range("L131").Value = 41913.332
$ws.Range("N131").Value = -51993.332
